$wb = $excel.ActiveWorkbook

# --- Sheet: Cases by Age Group ---
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value = 262
$ws.Range("B3").Value = 1251
$ws.Range("B4").Value = 3376
$ws.Range("B5").Value = 14683
$ws.Range("B6").Value = 16200
$ws.Range("B7").Value = 14164
$ws.Range("B8").Value = 11971
$ws.Range("B9").Value = 4320
$ws.Range("B10").Value = 2888
$ws.Range("B11").Value = 1714
$ws.Range("B12").Value = 1103
$ws.Range("B13").Value = 1716
$ws.Range("B20").Select()

# --- Sheet: Cases by Gender ---
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 24724
$ws.Range("B3").Value = 48028
$ws.Range("B4").Value = 910

# --- Sheet: Cases by RaceEthnicity ---
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 928
$ws.Range("B3").Value = 12369
$ws.Range("B4").Value = 27435
$ws.Range("B5").Value = 400
$ws.Range("B6").Value = 24275
$ws.Range("B7").Value = 8255
$ws.Range("B19").Select()

# --- Sheet: Fatalities by Age Group ---
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B4").Value = 26
$ws.Range("B5").Value = 192
$ws.Range("B6").Value = 632
$ws.Range("B7").Value = 1867
$ws.Range("B8").Value = 4320
$ws.Range("B9").Value = 3659
$ws.Range("B10").Value = 4689
$ws.Range("B11").Value = 5323
$ws.Range("B12").Value = 5351
$ws.Range("B13").Value = 14022
$ws.Range("C20").Select()

# --- Sheet: Fatalities by Gender ---
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 16854
$ws.Range("B3").Value = 23240
$ws.Range("D14").Select()

# --- Sheet: Fatalities by Race-Ethnicity ---
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 800
$ws.Range("B3").Value = 3816
$ws.Range("B4").Value = 18648
$ws.Range("B5").Value = 211
$ws.Range("B6").Value = 16598
$ws.Range("B7").Value = 22
$ws.Range("D17").Select()

# Make "Cases by Age Group" the active sheet with tabSelected, matching tabSelected attribute in the original file.
$wb.Worksheets.Item("Cases by Age Group").Activate()
$wb.Worksheets.Item("Cases by Age Group").Range("B20").Select()

$excel.CalculateFullRebuild()
